$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-18 Saturday" "2024-05-19 Sunday"

Replace-Text "826×4=3304" "956×5=4780"
Replace-Text "470×7=3290" "216×3=648"
Replace-Text "292×5=1460" "615×8=4920"
Replace-Text "657×8=5256" "356×4=1424"
Replace-Text "777×7=5439" "258×5=1290"

Replace-Text "253×9=2277" "171×3=513"
Replace-Text "731×8=5848" "340×5=1700"
Replace-Text "825×4=3300" "326×4=1304"
Replace-Text "186×3=558" "129×3=387"
Replace-Text "631×7=4417" "658×5=3290"

Replace-Text "341×8=2728" "593×3=1779"
Replace-Text "586×5=2930" "244×9=2196"
Replace-Text "181×3=543" "520×3=1560"
Replace-Text "209×6=1254" "952×3=2856"
Replace-Text "885×4=3540" "549×9=4941"

Replace-Text "115×2=230" "880×9=7920"
Replace-Text "555×5=2775" "168×8=1344"
Replace-Text "824×4=3296" "155×3=465"
Replace-Text "403×6=2418" "797×5=3985"
Replace-Text "651×7=4557" "791×9=7119"

Replace-Text "276×8=2208" "202×8=1616"
Replace-Text "772×7=5404" "464×4=1856"
Replace-Text "664×9=5976" "255×5=1275"
Replace-Text "323×5=1615" "919×8=7352"
Replace-Text "929×8=7432" "343×4=1372"
